$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A302").Value = 44376
$ws.Range("B302").Value = 1
$ws.Range("C302").Value = 29
$ws.Range("D302").Value = 15.34286001491961
$ws.Range("A303").Value = 44377
$ws.Range("B303").Value = 2
$ws.Range("C303").Value = 28
$ws.Range("D303").Value = 14.81379587647411
$ws.Range("A304").Value = 44378
$ws.Range("B304").Value = 3
$ws.Range("C304").Value = 29
$ws.Range("D304").Value = 15.34286001491961
$ws.Range("A305").Value = 44379
$ws.Range("B305").Value = 2
$ws.Range("C305").Value = 27
$ws.Range("D305").Value = 14.2847317380286
$ws.Range("A306").Value = 44380
$ws.Range("B306").Value = 1
$ws.Range("C306").Value = 18
$ws.Range("D306").Value = 9.523154492019067
$ws.Range("A307").Value = 44381
$ws.Range("B307").Value = 1
$ws.Range("C307").Value = 14
$ws.Range("D307").Value = 7.406897938237053
$ws.Range("A308").Value = 44382
$ws.Range("B308").Value = 0
$ws.Range("C308").Value = 10
$ws.Range("D308").Value = 5.290641384455038
$ws.Range("A309").Value = 44383
$ws.Range("B309").Value = 1
$ws.Range("C309").Value = 10
$ws.Range("D309").Value = 5.290641384455038
$ws.Range("A310").Value = 44384
$ws.Range("B310").Value = 5
$ws.Range("C310").Value = 13
$ws.Range("D310").Value = 6.877833799791548
$ws.Range("A311").Value = 44385
$ws.Range("B311").Value = 1
$ws.Range("C311").Value = 11
$ws.Range("D311").Value = 5.819705522900541
$ws.Range("A312").Value = 44386
$ws.Range("B312").Value = 3
$ws.Range("C312").Value = 12
$ws.Range("D312").Value = 6.348769661346044
$ws.Range("A313").Value = 44387
$ws.Range("B313").Value = 3
$ws.Range("C313").Value = 14
$ws.Range("D313").Value = 7.406897938237053
$ws.Range("A314").Value = 44388
$ws.Range("B314").Value = 7
$ws.Range("C314").Value = 20
$ws.Range("D314").Value = 10.58128276891008
$ws.Range("A315").Value = 44389
$ws.Range("B315").Value = 8
$ws.Range("C315").Value = 28
$ws.Range("D315").Value = 14.81379587647411
$ws.Range("A316").Value = 44390
$ws.Range("B316").Value = 3
$ws.Range("C316").Value = 30
$ws.Range("D316").Value = 15.87192415336511
$ws.Range("A317").Value = 44391
$ws.Range("B317").Value = 5
$ws.Range("C317").Value = 30
$ws.Range("D317").Value = 15.87192415336511
$ws.Range("A318").Value = 44392
$ws.Range("B318").Value = 5
$ws.Range("C318").Value = 34
$ws.Range("D318").Value = 17.98818070714713
$ws.Range("A319").Value = 44393
$ws.Range("B319").Value = 13
$ws.Range("C319").Value = 44
$ws.Range("D319").Value = 23.27882209160216
$ws.Range("A320").Value = 44394
$ws.Range("B320").Value = 9
$ws.Range("C320").Value = 50
$ws.Range("D320").Value = 26.45320692227519
$ws.Range("A321").Value = 44395
$ws.Range("B321").Value = 20
$ws.Range("C321").Value = 63
$ws.Range("D321").Value = 33.33104072206674
$ws.Range("A322").Value = 44396
$ws.Range("B322").Value = 7
$ws.Range("C322").Value = 62
$ws.Range("D322").Value = 32.80197658362123
$ws.Range("A323").Value = 44397
$ws.Range("B323").Value = 4
$ws.Range("C323").Value = 63
$ws.Range("D323").Value = 33.33104072206674
$ws.Range("A324").Value = 44398
$ws.Range("B324").Value = 6
$ws.Range("C324").Value = 64
$ws.Range("D324").Value = 33.86010486051224
$ws.Range("A325").Value = 44399
$ws.Range("B325").Value = 20
$ws.Range("C325").Value = 79
$ws.Range("D325").Value = 41.79606693719479
$ws.Range("A326").Value = 44400
$ws.Range("B326").Value = 16
$ws.Range("C326").Value = 82
$ws.Range("D326").Value = 43.38325935253131
$ws.Range("A327").Value = 44401
$ws.Range("B327").Value = 16
$ws.Range("C327").Value = 89
$ws.Range("D327").Value = 47.08670832164984
$ws.Range("A328").Value = 44402
$ws.Range("B328").Value = 22
$ws.Range("C328").Value = 91
$ws.Range("D328").Value = 48.14483659854084

$ws.Range("A301").Copy()
$ws.Range("A302:A328").PasteSpecial(-4122)
Write-Host "done"